$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu values per the diff (case with 380 kV done)
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028514446729083
$ws.Cells.Item(2, 4).Value = 1.037740055646572
$ws.Cells.Item(2, 5).Value = 1.028467046628688
$ws.Cells.Item(2, 6).Value = 1.044509931311164
$ws.Cells.Item(2, 9).Value = 1.03049677041457
$ws.Cells.Item(2, 10).Value = 1.033666341656593
$ws.Cells.Item(2, 11).Value = 1.040530078513809
$ws.Cells.Item(2, 12).Value = 1.031283779350623
$ws.Cells.Item(2, 13).Value = 1.0472807793711
$ws.Cells.Item(2, 14).Value = 1.015213766750359
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029446347452368
$ws.Cells.Item(3, 4).Value = 1.038641296716695
$ws.Cells.Item(3, 5).Value = 1.029257997927553
$ws.Cells.Item(3, 6).Value = 1.0455596232371
$ws.Cells.Item(3, 9).Value = 1.030586719395941
$ws.Cells.Item(3, 10).Value = 1.034238973821894
$ws.Cells.Item(3, 11).Value = 1.041240947626003
$ws.Cells.Item(3, 12).Value = 1.031882718062072
$ws.Cells.Item(3, 13).Value = 1.048141101815469
$ws.Cells.Item(3, 14).Value = 1.015405658170271
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03004972497751
$ws.Cells.Item(4, 4).Value = 1.039225176865902
$ws.Cells.Item(4, 5).Value = 1.029770499734586
$ws.Cells.Item(4, 6).Value = 1.046239891208616
$ws.Cells.Item(4, 9).Value = 1.030643452558911
$ws.Cells.Item(4, 10).Value = 1.034609272146156
$ws.Cells.Item(4, 11).Value = 1.041701012690289
$ws.Cells.Item(4, 12).Value = 1.032270320955704
$ws.Cells.Item(4, 13).Value = 1.048698233132253
$ws.Cells.Item(4, 14).Value = 1.015529685030047
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030303473811599
$ws.Cells.Item(5, 4).Value = 1.03947081046298
$ws.Cells.Item(5, 5).Value = 1.029986122553103
$ws.Cells.Item(5, 6).Value = 1.046526125395439
$ws.Cells.Item(5, 9).Value = 1.030666950770141
$ws.Cells.Item(5, 10).Value = 1.034764888843104
$ws.Cells.Item(5, 11).Value = 1.041894443395585
$ws.Cells.Item(5, 12).Value = 1.032433280039457
$ws.Cells.Item(5, 13).Value = 1.048932556628042
$ws.Cells.Item(5, 14).Value = 1.015581792100737
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030346084544051
$ws.Cells.Item(6, 4).Value = 1.039512063379827
$ws.Cells.Item(6, 5).Value = 1.030022336330475
$ws.Cells.Item(6, 6).Value = 1.046574199974278
$ws.Cells.Item(6, 9).Value = 1.030670875542361
$ws.Cells.Item(6, 10).Value = 1.034791014217857
$ws.Cells.Item(6, 11).Value = 1.041926922366525
$ws.Cells.Item(6, 12).Value = 1.032460642191868
$ws.Cells.Item(6, 13).Value = 1.048971906728609
$ws.Cells.Item(6, 14).Value = 1.015590539113086
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030053115233641
$ws.Cells.Item(7, 4).Value = 1.039228458366126
$ws.Cells.Item(7, 5).Value = 1.029773380242502
$ws.Cells.Item(7, 6).Value = 1.046243714904774
$ws.Cells.Item(7, 9).Value = 1.03064376792855
$ws.Cells.Item(7, 10).Value = 1.034611351728126
$ws.Cells.Item(7, 11).Value = 1.041703597246113
$ws.Cells.Item(7, 12).Value = 1.032272498381697
$ws.Cells.Item(7, 13).Value = 1.048701363761124
$ws.Cells.Item(7, 14).Value = 1.015530381420636
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.02882930862585
$ws.Cells.Item(8, 4).Value = 1.038044485093373
$ws.Cells.Item(8, 5).Value = 1.028734205866633
$ws.Cells.Item(8, 6).Value = 1.044864462111039
$ws.Cells.Item(8, 9).Value = 1.03052747311762
$ws.Cells.Item(8, 10).Value = 1.033859913049992
$ws.Cells.Item(8, 11).Value = 1.040770301783563
$ws.Cells.Item(8, 12).Value = 1.031486182628697
$ws.Cells.Item(8, 13).Value = 1.047571436541702
$ws.Cells.Item(8, 14).Value = 1.015278646060631
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026675722770029
$ws.Cells.Item(9, 4).Value = 1.035963718661752
$ws.Cells.Item(9, 5).Value = 1.026908494500353
$ws.Cells.Item(9, 6).Value = 1.04244211266595
$ws.Cells.Item(9, 9).Value = 1.030311312807202
$ws.Cells.Item(9, 10).Value = 1.032534039346529
$ws.Cells.Item(9, 11).Value = 1.039126414798383
$ws.Cells.Item(9, 12).Value = 1.030101017984645
$ws.Cells.Item(9, 13).Value = 1.045583817667021
$ws.Cells.Item(9, 14).Value = 1.014834002486229
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025242016266029
$ws.Cells.Item(10, 4).Value = 1.034580342940987
$ws.Cells.Item(10, 5).Value = 1.025695094819054
$ws.Cells.Item(10, 6).Value = 1.040832705982235
$ws.Cells.Item(10, 9).Value = 1.030159683900419
$ws.Cells.Item(10, 10).Value = 1.03164900714255
$ws.Cells.Item(10, 11).Value = 1.038031021462746
$ws.Cells.Item(10, 12).Value = 1.029177920205942
$ws.Cells.Item(10, 13).Value = 1.04426112626484
$ws.Cells.Item(10, 14).Value = 1.014536886829808
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024621695310981
$ws.Cells.Item(11, 4).Value = 1.033982242876078
$ws.Cells.Item(11, 5).Value = 1.025170583225026
$ws.Cells.Item(11, 6).Value = 1.040137131692903
$ws.Cells.Item(11, 9).Value = 1.030092249231479
$ws.Cells.Item(11, 10).Value = 1.031265525050257
$ws.Cells.Item(11, 11).Value = 1.037556842746133
$ws.Cells.Item(11, 12).Value = 1.02877830304457
$ws.Cells.Item(11, 13).Value = 1.043688965667505
$ws.Cells.Item(11, 14).Value = 1.014408074142852
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024391353970046
$ws.Cells.Item(12, 4).Value = 1.033760219896664
$ws.Cells.Item(12, 5).Value = 1.024975892559906
$ws.Cells.Item(12, 6).Value = 1.03987896221835
$ws.Cells.Item(12, 9).Value = 1.030066934305897
$ws.Cells.Item(12, 10).Value = 1.031123044954082
$ws.Cells.Item(12, 11).Value = 1.03738073274703
$ws.Cells.Item(12, 12).Value = 1.028629881976161
$ws.Cells.Item(12, 13).Value = 1.043476526915551
$ws.Cells.Item(12, 14).Value = 1.014360203743372
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024440759630658
$ws.Cells.Item(13, 4).Value = 1.03380783831491
$ws.Cells.Item(13, 5).Value = 1.025017648179833
$ws.Cells.Item(13, 6).Value = 1.039934331479677
$ws.Cells.Item(13, 9).Value = 1.030072376503644
$ws.Cells.Item(13, 10).Value = 1.031153609113443
$ws.Cells.Item(13, 11).Value = 1.037418507969534
$ws.Cells.Item(13, 12).Value = 1.028661718125596
$ws.Cells.Item(13, 13).Value = 1.043522091788207
$ws.Cells.Item(13, 14).Value = 1.0143704731703
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024602653719933
$ws.Cells.Item(14, 4).Value = 1.033963887565013
$ws.Cells.Item(14, 5).Value = 1.025154487245972
$ws.Cells.Item(14, 6).Value = 1.040115787282658
$ws.Cells.Item(14, 9).Value = 1.030090162129331
$ws.Cells.Item(14, 10).Value = 1.031253748366487
$ws.Cells.Item(14, 11).Value = 1.037542285002766
$ws.Cells.Item(14, 12).Value = 1.02876603421053
$ws.Cells.Item(14, 13).Value = 1.043671403629939
$ws.Cells.Item(14, 14).Value = 1.014404117640746
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024702411774121
$ws.Cells.Item(15, 4).Value = 1.03406005298209
$ws.Cells.Item(15, 5).Value = 1.025238816410404
$ws.Cells.Item(15, 6).Value = 1.040227614432152
$ws.Cells.Item(15, 9).Value = 1.030101085116368
$ws.Cells.Item(15, 10).Value = 1.031315442486557
$ws.Cells.Item(15, 11).Value = 1.037618550942038
$ws.Cells.Item(15, 12).Value = 1.028830308748998
$ws.Cells.Item(15, 13).Value = 1.043763411164473
$ws.Cells.Item(15, 14).Value = 1.01442484398846
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025283195182689
$ws.Cells.Item(16, 4).Value = 1.034620056162102
$ws.Cells.Item(16, 5).Value = 1.025729923967709
$ws.Cells.Item(16, 6).Value = 1.040878896635535
$ws.Cells.Item(16, 9).Value = 1.03016412189132
$ws.Cells.Item(16, 10).Value = 1.031674452240397
$ws.Cells.Item(16, 11).Value = 1.038062494045953
$ws.Cells.Item(16, 12).Value = 1.029204443466035
$ws.Cells.Item(16, 13).Value = 1.044299110848311
$ws.Cells.Item(16, 14).Value = 1.014545432370164
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025647635297492
$ws.Cells.Item(17, 4).Value = 1.034971576031985
$ws.Cells.Item(17, 5).Value = 1.026038224235251
$ws.Cells.Item(17, 6).Value = 1.041287780073097
$ws.Cells.Item(17, 9).Value = 1.030203187425615
$ws.Cells.Item(17, 10).Value = 1.031899581310102
$ws.Cells.Item(17, 11).Value = 1.038341004349006
$ws.Cells.Item(17, 12).Value = 1.0294391531463
$ws.Cells.Item(17, 13).Value = 1.044635295361036
$ws.Cells.Item(17, 14).Value = 1.014621031783343
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025860253636443
$ws.Cells.Item(18, 4).Value = 1.035176699419575
$ws.Cells.Item(18, 5).Value = 1.026218137167371
$ws.Cells.Item(18, 6).Value = 1.041526401133697
$ws.Cells.Item(18, 9).Value = 1.030225802072797
$ws.Cells.Item(18, 10).Value = 1.032030870401732
$ws.Cells.Item(18, 11).Value = 1.038503467543039
$ws.Cells.Item(18, 12).Value = 1.029576063975925
$ws.Cells.Item(18, 13).Value = 1.04483144136713
$ws.Cells.Item(18, 14).Value = 1.014665112214725
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025932758871351
$ws.Cells.Item(19, 4).Value = 1.035246656049354
$ws.Cells.Item(19, 5).Value = 1.026279497477234
$ws.Cells.Item(19, 6).Value = 1.041607786134075
$ws.Cells.Item(19, 9).Value = 1.030233483958762
$ws.Cells.Item(19, 10).Value = 1.032075632341178
$ws.Cells.Item(19, 11).Value = 1.038558865442839
$ws.Cells.Item(19, 12).Value = 1.029622748467216
$ws.Cells.Item(19, 13).Value = 1.044898331430111
$ws.Cells.Item(19, 14).Value = 1.014680139875207
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025608529497406
$ws.Cells.Item(20, 4).Value = 1.034933852187085
$ws.Cells.Item(20, 5).Value = 1.026005137567132
$ws.Cells.Item(20, 6).Value = 1.04124389768156
$ws.Cells.Item(20, 9).Value = 1.030199013811208
$ws.Cells.Item(20, 10).Value = 1.031875429659164
$ws.Cells.Item(20, 11).Value = 1.038311121493402
$ws.Cells.Item(20, 12).Value = 1.029413970117703
$ws.Cells.Item(20, 13).Value = 1.044599220220547
$ws.Cells.Item(20, 14).Value = 1.014612922274934
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024554977910425
$ws.Cells.Item(21, 4).Value = 1.033917931130867
$ws.Cells.Item(21, 5).Value = 1.025114187778707
$ws.Cells.Item(21, 6).Value = 1.040062347602924
$ws.Cells.Item(21, 9).Value = 1.03008493206379
$ws.Cells.Item(21, 10).Value = 1.031224260887967
$ws.Cells.Item(21, 11).Value = 1.037505835161477
$ws.Cells.Item(21, 12).Value = 1.028735315314358
$ws.Cells.Item(21, 13).Value = 1.043627432599655
$ws.Cells.Item(21, 14).Value = 1.014394210832295
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023892993779286
$ws.Cells.Item(22, 4).Value = 1.03327998034149
$ws.Cells.Item(22, 5).Value = 1.024554801492041
$ws.Cells.Item(22, 6).Value = 1.03932060503485
$ws.Cells.Item(22, 9).Value = 1.030011661615347
$ws.Cells.Item(22, 10).Value = 1.030814627135006
$ws.Cells.Item(22, 11).Value = 1.036999642185377
$ws.Cells.Item(22, 12).Value = 1.028308702701194
$ws.Cells.Item(22, 13).Value = 1.043016935987242
$ws.Cells.Item(22, 14).Value = 1.014256561713878
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024243883467668
$ws.Cells.Item(23, 4).Value = 1.033618093983828
$ws.Cells.Item(23, 5).Value = 1.02485126750195
$ws.Cells.Item(23, 6).Value = 1.039713707929667
$ws.Cells.Item(23, 9).Value = 1.030050649735237
$ws.Cells.Item(23, 10).Value = 1.0310318021158
$ws.Cells.Item(23, 11).Value = 1.03726797271903
$ws.Cells.Item(23, 12).Value = 1.02853484977429
$ws.Cells.Item(23, 13).Value = 1.043340523603286
$ws.Cells.Item(23, 14).Value = 1.014329544941869
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025626199581268
$ws.Cells.Item(24, 4).Value = 1.034950897697618
$ws.Cells.Item(24, 5).Value = 1.026020087740801
$ws.Cells.Item(24, 6).Value = 1.041263725855315
$ws.Cells.Item(24, 9).Value = 1.030200900218223
$ws.Cells.Item(24, 10).Value = 1.031886342827589
$ws.Cells.Item(24, 11).Value = 1.038324624229938
$ws.Cells.Item(24, 12).Value = 1.029425349217712
$ws.Cells.Item(24, 13).Value = 1.044615520852996
$ws.Cells.Item(24, 14).Value = 1.014616586660413
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02723212484348
$ws.Cells.Item(25, 4).Value = 1.036500981916222
$ws.Cells.Item(25, 5).Value = 1.027379831028951
$ws.Cells.Item(25, 6).Value = 1.043067384808477
$ws.Cells.Item(25, 9).Value = 1.030368523513698
$ws.Cells.Item(25, 10).Value = 1.032877010305226
$ws.Cells.Item(25, 11).Value = 1.039551309841276
$ws.Cells.Item(25, 12).Value = 1.030459059784129
$ws.Cells.Item(25, 13).Value = 1.046097248563448
$ws.Cells.Item(25, 14).Value = 1.014949076125733
